$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E25").Value = 23.36
$ws.Range("E26").Value = 57.7

$cos = $ws.ChartObjects()
$co = $cos.Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
$s = $sc.Item(4)

$arr = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,6,0,0,13.8,30.56,17.5,67.86,23.36,57.7)
$s.Values = $arr
Write-Host ("Formula after: " + $s.Formula)
